# This script applies the edits described by the commit diff to the
# "tables.xlsx" style workbook (3 sheets: Table 1, Table 2, Table 3).
#
# Sheet 1 ("Table 1"):
#   - rename some row labels
#   - for several binary yes/no variables, drop the "no: .." line and the
#     "yes: " prefix, keeping only the percentage line for "yes"
#
# Sheet 2 ("Table 2"):
#   - rename several "# ..." abbreviated labels to full text
#   - same "no/yes" simplification for two rows
#
# Sheet 3 ("Table 3"):
#   - two brand new rows are inserted (categorical breakdowns for
#     "Overall Mental Health" and "Quality of Life")
#   - remaining rows are renamed / renumbered, and the yes/no rows are
#     simplified to show only the "yes" percentage
#   - a couple of p-values are updated to new numbers

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Table 1
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Table 1")

# Row 8: "Sum of co-morbidities" -> "Number of co-morbidities"
$ws1.Range("A8").Value = "Number of co-morbidities"

# Row 10: "Pre-CoV depression/anxiety" -> "Depression/anxiety before COVID-19"
$ws1.Range("A10").Value = "Depression/anxiety before COVID-19"
$ws1.Range("B10").Value = "5.96% (69) `nn = 1157"
$ws1.Range("C10").Value = "4.59% (41) `nn = 893"

# Row 11: "Pre-CoV sleep disorders" -> "Sleep disorders before COVID-19"
$ws1.Range("A11").Value = "Sleep disorders before COVID-19"
$ws1.Range("B11").Value = "4.58% (53) `nn = 1157"
$ws1.Range("C11").Value = "4.03% (36) `nn = 893"

# Row 12: Bruxism - simplify to "yes" only
$ws1.Range("B12").Value = "7.17% (83) `nn = 1157"
$ws1.Range("C12").Value = "5.26% (47) `nn = 893"

# Row 14: Hypertension - simplify to "yes" only
$ws1.Range("B14").Value = "11.2% (130) `nn = 1157"
$ws1.Range("C14").Value = "9.41% (84) `nn = 893"

# Row 15: Cardiovascular disease - simplify to "yes" only
$ws1.Range("B15").Value = "2.94% (34) `nn = 1157"
$ws1.Range("C15").Value = "2.91% (26) `nn = 893"

# Row 16: Pulmonary disease - simplify to "yes" only
$ws1.Range("B16").Value = "4.15% (48) `nn = 1157"
$ws1.Range("C16").Value = "2.58% (23) `nn = 893"

# Row 17: Hay fever/allergy - simplify to "yes" only
$ws1.Range("B17").Value = "18% (208) `nn = 1157"
$ws1.Range("C17").Value = "11.4% (102) `nn = 893"

foreach ($r in 10,11,12,14,15,16,17) {
    $ws1.Rows.Item($r).AutoFit()
}

# ---------------------------------------------------------------------
# Table 2
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Table 2")

# Row 2: Acute COVID-19 symptoms - simplify to "yes" only
$ws2.Range("B2").Value = "91.7% (1060) `nn = 1156"
$ws2.Range("C2").Value = "87.7% (782) `nn = 892"

# Row 3: "# acute symptoms" -> "Number of acute symptoms"
$ws2.Range("A3").Value = "Number of acute symptoms"

# Row 4: "# acute NC" -> "Number of acute neurocognitive symptoms"
$ws2.Range("A4").Value = "Number of acute neurocognitive symptoms"

# Row 5: "# acute NC" -> "Number of acute neurocognitive symptoms"
$ws2.Range("A5").Value = "Number of acute neurocognitive symptoms"

# Row 6: Persistent COVID-19 symptoms - simplify to "yes" only
$ws2.Range("B6").Value = "47.6% (550) `nn = 1156"
$ws2.Range("C6").Value = "49.3% (440) `nn = 892"

# Row 7: "# persistent symptoms" -> "Number of persistent symptoms"
$ws2.Range("A7").Value = "Number of persistent symptoms"

# Row 8: "# persist. NC" -> "Number of persistent neurocognitive symptoms"
$ws2.Range("A8").Value = "Number of persistent neurocognitive symptoms"

# Row 9: "# persist. NC" -> "Number of persistent neurocognitive symptoms"
$ws2.Range("A9").Value = "Number of persistent neurocognitive symptoms"

foreach ($r in 2,6) {
    $ws2.Rows.Item($r).AutoFit()
}

# ---------------------------------------------------------------------
# Table 3
# ---------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("Table 3")

# Insert a fresh row above the old row 2 ("OMH score") for the new
# "Overall Mental Health" categorical breakdown, and a fresh row above
# the old row 3 ("QoL score", now shifted to row 4) for the new
# "Quality of Life" categorical breakdown.
$ws3.Rows.Item(2).Insert()
$ws3.Rows.Item(4).Insert()

# New row 2: Overall Mental Health (categorical)
$ws3.Range("A2:E2").Style = "Normal"
$ws3.Range("A2").Value = "Overall Mental Health"
$ws3.Range("B2").Value = "poor: 3.46% (40)`nfair: 18.3% (212)`ngood: 48.6% (562)`nexcellent: 29.6% (343) `nn = 1157"
$ws3.Range("C2").Value = "poor: 2.91% (26)`nfair: 21.2% (189)`ngood: 48.2% (430)`nexcellent: 27.8% (248) `nn = 893"
$ws3.Range("D2").Value = "Chi"
$ws3.Range("E2").Value = "ns"

# Row 3 (was row 2): OMH score -> Overall Mental Health Score
$ws3.Range("A3").Value = "Overall Mental Health Score"

# New row 4: Quality of Life (categorical)
$ws3.Range("A4:E4").Style = "Normal"
$ws3.Range("A4").Value = "Quality of Life"
$ws3.Range("B4").Value = "poor: 4.32% (50)`nfair: 16% (185)`ngood: 51% (590)`nexcellent: 28.7% (332) `nn = 1157"
$ws3.Range("C4").Value = "poor: 3.36% (30)`nfair: 22.5% (201)`ngood: 54.3% (485)`nexcellent: 19.8% (177) `nn = 893"
$ws3.Range("D4").Value = "Chi"
$ws3.Range("E4").Value = "p = 8.3e-06"

# Row 5 (was row 3): QoL score -> Quality of Life Score; p-value updated
$ws3.Range("A5").Value = "Quality of Life Score"
$ws3.Range("E5").Value = "p = 2.1e-05"

# Row 6 (was row 4): DPR score - unchanged name; p-value updated
$ws3.Range("E6").Value = "p = 0.0076"

# Row 7 (was row 5): DPR+ -> Depression Screening-positive, simplified
$ws3.Range("A7").Value = "Depression Screening-positive"
$ws3.Range("B7").Value = "17.3% (200) `nn = 1154"
$ws3.Range("C7").Value = "23.2% (207) `nn = 892"

# Row 8 (was row 6): ANX score -> Anxiety score; p-value updated
$ws3.Range("A8").Value = "Anxiety score"
$ws3.Range("E8").Value = "p = 5.2e-09"

# Row 9 (was row 7): ANX+ -> Anxiety Screening-positive, simplified; p-value updated
$ws3.Range("A9").Value = "Anxiety Screening-positive"
$ws3.Range("B9").Value = "12.4% (143) `nn = 1151"
$ws3.Range("C9").Value = "19.3% (172) `nn = 893"
$ws3.Range("E9").Value = "p = 7.1e-05"

# Row 10 (was row 8): Stress score -> Psychosocial Stress Score
$ws3.Range("A10").Value = "Psychosocial Stress Score"

# Row 11 (was row 9): Substantial psychosocial stress, simplified; p-value updated
$ws3.Range("B11").Value = "21.3% (246) `nn = 1153"
$ws3.Range("C11").Value = "25.6% (228) `nn = 890"
$ws3.Range("E11").Value = "p = 0.038"

foreach ($r in 2,3,4,5,6,7,8,9,10,11) {
    $ws3.Rows.Item($r).AutoFit()
}
